{"js": "// Remove the page-break paragraph (and the empty paragraph that follows\n// it) that sits between the \"[AcademySection]\" placeholder paragraph and\n// the \"Financial assessment\" Heading 2 paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Find the paragraph that consists solely of a manual page break\n// (Word/Office.js represents a page-break run as \"\\f\" in paragraph.text)\n// immediately followed by the \"Financial assessment\" heading, which is\n// the pattern this specific edit targets (there are other, unrelated\n// page breaks earlier in the document that must stay untouched).\nlet targetIndex = -1;\nfor (let i = 0; i < items.length - 1; i++) {\n  if (items[i].text === \"\\f\" && items[i + 1].text === \"\") {\n    const after = items[i + 2];\n    if (after) {\n      after.load(\"text\");\n    }\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < items.length - 2; i++) {\n  if (\n    items[i].text === \"\\f\" &&\n    items[i + 1].text === \"\" &&\n    items[i + 2] &&\n    items[i + 2].text === \"Financial assessment\"\n  ) {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the page-break paragraph to remove.\");\n}\n\n// Delete the page-break paragraph itself, then the now-adjacent empty\n// paragraph that originally followed it.\nitems[targetIndex].delete();\nitems[targetIndex + 1].delete();\nawait context.sync();\n", "ps1": "# Remove the page-break paragraph (and the empty paragraph that follows\n# it) that sits between the \"[AcademySection]\" placeholder paragraph and\n# the \"Financial assessment\" Heading 2 paragraph.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$targetIndex = -1\n\n# Locate the paragraph that consists solely of a manual page break\n# (a page-break run shows up as a form-feed, chr(12), in Range.Text)\n# that is immediately followed by an otherwise-empty paragraph and then\n# the \"Financial assessment\" heading. There are other, unrelated page\n# breaks earlier in the document that must stay untouched, so we key off\n# this specific surrounding context rather than the first page break we\n# find.\nfor ($i = 1; $i -le $count - 2; $i++) {\n    $t0 = $d.Paragraphs.Item($i).Range.Text\n    $t1 = $d.Paragraphs.Item($i + 1).Range.Text\n    $t2 = $d.Paragraphs.Item($i + 2).Range.Text\n\n    if ($t0.Contains([char]12) -and `\n        $t1.Trim() -eq \"\" -and -not $t1.Contains([char]12) -and `\n        $t2.Trim() -eq \"Financial assessment\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the page-break paragraph to remove.\"\n}\n\n# Delete the trailing empty paragraph first so the page-break paragraph's\n# index doesn't shift, then delete the page-break paragraph itself.\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n$d.Paragraphs.Item($targetIndex).Range.Delete()\n"}
